$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.532203912734985
$ws.Range("B1").Value = 1.738739490509033
$ws.Range("C1").Value = 1.517043948173523
$ws.Range("D1").Value = 1.432864904403687
$ws.Range("E1").Value = 0.8173864483833313
